# Add the new "why plant-based? / safety" translation rows to the sheet.
# Each new row duplicates the formatting of the last existing row (B72, style s="3")
# by copying that row and inserting it, then overwriting the value with the new text.
# This keeps formatting/styles identical to the author's edit while adding fresh
# shared-string entries for each new piece of text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$texts = @(
  "Conventional meat-based pet food uses ingredients not suited for human consumption. This means that there are much less regulations in place which results in <b>wide-spread contamination</b> with toxic compounds such as:",
  "An ongoing and increasing contamination with toxic elements. <b>Heavy metals</b> including arsenic, cadmium, nickel, lead, antimony, <b>radioactive elements</b> such as uranium and thorium and other elements such as beryllium routinely go way above safe limits. [1,2]",
  "Endocrine disrupting compounds such as PCBs and PBDEs. [5]",
  "Nitrate and nitrite used in meat processing are routinely <b>two to three times the safety limit</b>. [2] Chronic exposure can result in cardiac and thyroid diseases as well as cancer. [4]",
  "Mycotoxin contamination</b> above safe limits of all products irrespective of marketing channels. The long-term exposure to mycotoxins is implicated in numerous clinical conditions such as vomiting, reduced immunity and cancer. [3]",
  "VeggieAnimals plant-based pet food is formulated from human grade ingredients which ensures that your pet is kept as healthy as possible and does not suffer needlessly from contaminants.",
  "References:"
)

$lastRow = 72
$row = $lastRow + 1
foreach ($t in $texts) {
  $ws.Rows($lastRow).Copy()
  $ws.Rows($row).Insert(-4121)
  $ws.Range("B$row").Value = $t
  $row = $row + 1
}

$lastNewRow = $row - 1

# Update the view so the newly added content is visible, matching the saved
# workbook state (scrolled down, with the cell after the last new row selected).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 40
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B" + ($lastNewRow + 1)).Select()
